# Update cryptocurrency price/volume data per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.014.48'
$ws.Range('E2').Value = '  +4.31%  '
$ws.Range('D3').Value = '3.251.97'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '578.52'
$ws.Range('E5').Value = '  +3.02%  '
$ws.Range('D6').Value = '176.96'
$ws.Range('E6').Value = '  +2.49%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.604'
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').Value = '3.248.16'
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('E10').Value = '  +4.00%  '
$ws.Range('D11').Value = '6.72'
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('E12').Value = '  +2.59%  '
$ws.Range('D13').Value = '3.819.63'
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('D14').Value = '0.136'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').Value = '27.91'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '66.977.88'
$ws.Range('E16').Value = '  +4.27%  '
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('D18').Value = '3.252.54'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.80'
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').Value = '13.35'
$ws.Range('E20').Value = '  +1.81%  '
$ws.Range('D21').Value = '369.33'
$ws.Range('E21').Value = '  +4.46%  '
$ws.Range('D22').Value = '7.51'
$ws.Range('E22').Value = '  +4.48%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '70.66'
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '0.507'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.388.13'
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.80'
$ws.Range('E28').Value = '  +2.93%  '
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +4.51%  '
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('D33').Value = '22.48'
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '174.06'
$ws.Range('E35').Value = '  +10.77%  '
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('D37').Value = '6.76'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('E38').Value = '  +4.69%  '
$ws.Range('D39').Value = '0.854'
$ws.Range('E39').Value = '  +6.84%  '
$ws.Range('E40').Value = '  +9.34%  '
$ws.Range('E41').Value = '  +2.61%  '
$ws.Range('D42').Value = '2.57'
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.730.31'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '6.43'
$ws.Range('E44').Value = '  +6.76%  '
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('D46').Value = '40.42'
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('E47').Value = '  +3.10%  '
$ws.Range('D48').Value = '24.66'
$ws.Range('E48').Value = '  +3.53%  '
$ws.Range('D49').Value = '334.31'
$ws.Range('E49').Value = '  +1.56%  '
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('E51').Value = '  +2.24%  '
